{"js": "// Office.js (Word JavaScript API) script.\n// Replaces the date line and the 25 \"A\u00d7B=C\" multiplication-table answers\n// with their updated values, per the commit diff.\n//\n// Every old value below is unique within the document, so a plain\n// case-sensitive whole-document search/replace for each pair is\n// unambiguous and safe.\n\nconst replacements = [\n  [\"2025-10-29 Wednesday\", \"2025-10-30 Thursday\"],\n  [\"111\u00d79=999\", \"171\u00d72=342\"],\n  [\"458\u00d74=1832\", \"151\u00d78=1208\"],\n  [\"941\u00d73=2823\", \"981\u00d73=2943\"],\n  [\"572\u00d74=2288\", \"678\u00d79=6102\"],\n  [\"242\u00d79=2178\", \"245\u00d77=1715\"],\n  [\"814\u00d79=7326\", \"465\u00d75=2325\"],\n  [\"353\u00d78=2824\", \"848\u00d76=5088\"],\n  [\"169\u00d72=338\", \"602\u00d73=1806\"],\n  [\"668\u00d79=6012\", \"581\u00d74=2324\"],\n  [\"795\u00d79=7155\", \"365\u00d77=2555\"],\n  [\"433\u00d76=2598\", \"709\u00d74=2836\"],\n  [\"841\u00d75=4205\", \"299\u00d76=1794\"],\n  [\"834\u00d72=1668\", \"525\u00d79=4725\"],\n  [\"118\u00d74=472\", \"340\u00d77=2380\"],\n  [\"543\u00d77=3801\", \"699\u00d78=5592\"],\n  [\"539\u00d79=4851\", \"554\u00d75=2770\"],\n  [\"802\u00d78=6416\", \"507\u00d78=4056\"],\n  [\"145\u00d74=580\", \"693\u00d79=6237\"],\n  [\"256\u00d78=2048\", \"602\u00d73=1806\"],\n  [\"653\u00d74=2612\", \"443\u00d79=3987\"],\n  [\"581\u00d73=1743\", \"837\u00d74=3348\"],\n  [\"390\u00d73=1170\", \"720\u00d78=5760\"],\n  [\"970\u00d78=7760\", \"788\u00d73=2364\"],\n  [\"753\u00d75=3765\", \"414\u00d73=1242\"],\n  [\"267\u00d72=534\", \"231\u00d74=924\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Replaces the date line and the 25 \"A\u00d7B=C\" multiplication-table answers\n# with their updated values, per the commit diff.\n#\n# Every old value below is unique within the document, so a plain\n# case-sensitive whole-document Find/Replace for each pair is\n# unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"2025-10-29 Wednesday\", \"2025-10-30 Thursday\"),\n    @(\"111\u00d79=999\", \"171\u00d72=342\"),\n    @(\"458\u00d74=1832\", \"151\u00d78=1208\"),\n    @(\"941\u00d73=2823\", \"981\u00d73=2943\"),\n    @(\"572\u00d74=2288\", \"678\u00d79=6102\"),\n    @(\"242\u00d79=2178\", \"245\u00d77=1715\"),\n    @(\"814\u00d79=7326\", \"465\u00d75=2325\"),\n    @(\"353\u00d78=2824\", \"848\u00d76=5088\"),\n    @(\"169\u00d72=338\", \"602\u00d73=1806\"),\n    @(\"668\u00d79=6012\", \"581\u00d74=2324\"),\n    @(\"795\u00d79=7155\", \"365\u00d77=2555\"),\n    @(\"433\u00d76=2598\", \"709\u00d74=2836\"),\n    @(\"841\u00d75=4205\", \"299\u00d76=1794\"),\n    @(\"834\u00d72=1668\", \"525\u00d79=4725\"),\n    @(\"118\u00d74=472\", \"340\u00d77=2380\"),\n    @(\"543\u00d77=3801\", \"699\u00d78=5592\"),\n    @(\"539\u00d79=4851\", \"554\u00d75=2770\"),\n    @(\"802\u00d78=6416\", \"507\u00d78=4056\"),\n    @(\"145\u00d74=580\", \"693\u00d79=6237\"),\n    @(\"256\u00d78=2048\", \"602\u00d73=1806\"),\n    @(\"653\u00d74=2612\", \"443\u00d79=3987\"),\n    @(\"581\u00d73=1743\", \"837\u00d74=3348\"),\n    @(\"390\u00d73=1170\", \"720\u00d78=5760\"),\n    @(\"970\u00d78=7760\", \"788\u00d73=2364\"),\n    @(\"753\u00d75=3765\", \"414\u00d73=1242\"),\n    @(\"267\u00d72=534\", \"231\u00d74=924\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
